$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 299, shifting existing rows 299:386 down to 300:387
$ws.Rows.Item(299).Insert()

# Populate the new row 299 with the new record's data
$ws.Cells.Item(299, 1).Value = 7
$ws.Cells.Item(299, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(299, 3).Value = "Ñuble"
$ws.Cells.Item(299, 4).Value = 45204
$ws.Cells.Item(299, 5).Value = 16
$ws.Cells.Item(299, 6).Value = 100112045
$ws.Cells.Item(299, 7).Value = "Zapallo"
$ws.Cells.Item(299, 8).Value = "Paine"
$ws.Cells.Item(299, 9).Value = "1a (guarda)"
$ws.Cells.Item(299, 10).Value = 500
$ws.Cells.Item(299, 11).Value = 450
$ws.Cells.Item(299, 12).Value = 450
$ws.Cells.Item(299, 13).Value = 450
$ws.Cells.Item(299, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(299, 15).Value = "Región del Maule"
$ws.Cells.Item(299, 16).Value = 450
$ws.Cells.Item(299, 17).Value = 1
$ws.Cells.Item(299, 18).Value = "Hortaliza"

Write-Host "Done"
